$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels: G1 was "latitude", now "longitude"; H1 was "longitude", now "latitude"
$ws.Range("G1").Value = "longitude"
$ws.Range("H1").Value = "latitude"

# Update G/H coordinate values for rows 2-101 (columns swapped meaning + higher precision values)
$ws.Range("G2").Value = 5.3586322220931768
$ws.Range("H2").Value = 51.61987376615707
$ws.Range("G3").Value = 5.3762496256572012
$ws.Range("H3").Value = 51.615141649572799
$ws.Range("G4").Value = 5.4823233183292182
$ws.Range("H4").Value = 51.646459266511897
$ws.Range("G5").Value = 5.4630797398896354
$ws.Range("H5").Value = 51.652550619523979
$ws.Range("G6").Value = 5.4859735903765889
$ws.Range("H6").Value = 51.649369557643517
$ws.Range("G7").Value = 5.4741942721130652
$ws.Range("H7").Value = 51.653823267353928
$ws.Range("G8").Value = 5.4718832989536104
$ws.Range("H8").Value = 51.651594611510419
$ws.Range("G9").Value = 5.4788942578157753
$ws.Range("H9").Value = 51.650441344485401
$ws.Range("G10").Value = 5.4682762333412942
$ws.Range("H10").Value = 51.649079970519907
$ws.Range("G11").Value = 5.4697950068818928
$ws.Range("H11").Value = 51.656876273651577
$ws.Range("G12").Value = 5.4932503938580313
$ws.Range("H12").Value = 51.647873868199667
$ws.Range("G13").Value = 5.4804119344502586
$ws.Range("H13").Value = 51.654802277044404
$ws.Range("G14").Value = 5.4852072999730463
$ws.Range("H14").Value = 51.653162558037572
$ws.Range("G15").Value = 5.4990592516330636
$ws.Range("H15").Value = 51.638676281886113
$ws.Range("G16").Value = 5.4512774659772578
$ws.Range("H16").Value = 51.655385590824643
$ws.Range("G17").Value = 5.4702828369908989
$ws.Range("H17").Value = 51.676417242106233
$ws.Range("G18").Value = 5.5260235450763426
$ws.Range("H18").Value = 51.648637495970263
$ws.Range("G19").Value = 5.5132650434779453
$ws.Range("H19").Value = 51.664265301658617
$ws.Range("G20").Value = 5.5080405362414346
$ws.Range("H20").Value = 51.678795079792692
$ws.Range("G21").Value = 5.5781551047767994
$ws.Range("H21").Value = 51.573853220935973
$ws.Range("G22").Value = 5.4241727275527856
$ws.Range("H22").Value = 51.629297799746197
$ws.Range("G23").Value = 5.4190194549987876
$ws.Range("H23").Value = 51.624945391976759
$ws.Range("G24").Value = 5.4079437083828319
$ws.Range("H24").Value = 51.612695266653759
$ws.Range("G25").Value = 5.473927301409101
$ws.Range("H25").Value = 51.62430361985033
$ws.Range("G26").Value = 5.4136497874507441
$ws.Range("H26").Value = 51.642040515163337
$ws.Range("G27").Value = 5.3966596001118399
$ws.Range("H27").Value = 51.628219072108308
$ws.Range("G28").Value = 5.4243722837192667
$ws.Range("H28").Value = 51.601090739414254
$ws.Range("G29").Value = 5.4476153591554413
$ws.Range("H29").Value = 51.637903051954922
$ws.Range("G30").Value = 5.4731844089914388
$ws.Range("H30").Value = 51.601310651299492
$ws.Range("G31").Value = 5.4447133137722119
$ws.Range("H31").Value = 51.619929405567127
$ws.Range("G32").Value = 5.4436335465204069
$ws.Range("H32").Value = 51.616330822455758
$ws.Range("G33").Value = 5.4271185150327073
$ws.Range("H33").Value = 51.617136043679608
$ws.Range("G34").Value = 5.4245833293788452
$ws.Range("H34").Value = 51.623457339020739
$ws.Range("G35").Value = 5.4667230656342456
$ws.Range("H35").Value = 51.610971795137829
$ws.Range("G36").Value = 5.4685519522894417
$ws.Range("H36").Value = 51.606474803882207
$ws.Range("G37").Value = 5.414022494553425
$ws.Range("H37").Value = 51.627093836939977
$ws.Range("G38").Value = 5.4512050535109964
$ws.Range("H38").Value = 51.614132970024137
$ws.Range("G39").Value = 5.4311553548430496
$ws.Range("H39").Value = 51.625550952105577
$ws.Range("G40").Value = 5.452461707508812
$ws.Range("H40").Value = 51.610101402095587
$ws.Range("G41").Value = 5.417220283356535
$ws.Range("H41").Value = 51.632610603721062
$ws.Range("G42").Value = 5.4414771793891141
$ws.Range("H42").Value = 51.613275387608788
$ws.Range("G43").Value = 5.4436208412954548
$ws.Range("H43").Value = 51.609944031411509
$ws.Range("G44").Value = 5.4362925389484094
$ws.Range("H44").Value = 51.614187827887037
$ws.Range("G45").Value = 5.4437974343889186
$ws.Range("H45").Value = 51.605980899384207
$ws.Range("G46").Value = 5.4369897906255478
$ws.Range("H46").Value = 51.618894428623094
$ws.Range("G47").Value = 5.4330788253202602
$ws.Range("H47").Value = 51.617556581600688
$ws.Range("G48").Value = 5.4334222598473696
$ws.Range("H48").Value = 51.622575800659853
$ws.Range("G49").Value = 5.4285090798288103
$ws.Range("H49").Value = 51.62065078716887
$ws.Range("G50").Value = 5.4674400005127346
$ws.Range("H50").Value = 51.617285144630458
$ws.Range("G51").Value = 5.4607067622473542
$ws.Range("H51").Value = 51.569190875717737
$ws.Range("G52").Value = 5.419280739417708
$ws.Range("H52").Value = 51.550543947786522
$ws.Range("G53").Value = 5.4281173323854244
$ws.Range("H53").Value = 51.551966484480587
$ws.Range("G54").Value = 5.5156764659981308
$ws.Range("H54").Value = 51.570692905258568
$ws.Range("G55").Value = 5.467137048710951
$ws.Range("H55").Value = 51.547821923875212
$ws.Range("G56").Value = 5.4093159781243481
$ws.Range("H56").Value = 51.584285732738699
$ws.Range("G57").Value = 5.4437599594075001
$ws.Range("H57").Value = 51.580559663505298
$ws.Range("G58").Value = 5.473587605777583
$ws.Range("H58").Value = 51.586997652594548
$ws.Range("G59").Value = 5.5102438079310776
$ws.Range("H59").Value = 51.543927264682132
$ws.Range("G60").Value = 5.4491157209396874
$ws.Range("H60").Value = 51.569542374889558
$ws.Range("G61").Value = 5.4705129147903326
$ws.Range("H61").Value = 51.568612896162222
$ws.Range("G62").Value = 5.468521297268051
$ws.Range("H62").Value = 51.574044649912338
$ws.Range("G63").Value = 5.4525935036505686
$ws.Range("H63").Value = 51.5764625006736
$ws.Range("G64").Value = 5.4570175681632112
$ws.Range("H64").Value = 51.573102752310092
$ws.Range("G65").Value = 5.4821666683244983
$ws.Range("H65").Value = 51.55279536850837
$ws.Range("G66").Value = 5.4070103184239304
$ws.Range("H66").Value = 51.581207786284587
$ws.Range("G67").Value = 5.4592693853583674
$ws.Range("H67").Value = 51.565071866869431
$ws.Range("G68").Value = 5.4621773056410889
$ws.Range("H68").Value = 51.559756599785743
$ws.Range("G69").Value = 5.5624257850291752
$ws.Range("H69").Value = 51.634017535585798
$ws.Range("G70").Value = 5.5522559207141002
$ws.Range("H70").Value = 51.631636131591179
$ws.Range("G71").Value = 5.5186563781892497
$ws.Range("H71").Value = 51.62056598325664
$ws.Range("G72").Value = 5.5317124325831042
$ws.Range("H72").Value = 51.612109493221119
$ws.Range("G73").Value = 5.5615688762907283
$ws.Range("H73").Value = 51.641636838131802
$ws.Range("G74").Value = 5.4994588238949156
$ws.Range("H74").Value = 51.600394256069677
$ws.Range("G75").Value = 5.5693618094927224
$ws.Range("H75").Value = 51.605964287604777
$ws.Range("G76").Value = 5.5566587337355751
$ws.Range("H76").Value = 51.56770418353598
$ws.Range("G77").Value = 5.5831975436145527
$ws.Range("H77").Value = 51.629000356055229
$ws.Range("G78").Value = 5.5248305073567181
$ws.Range("H78").Value = 51.628239100171378
$ws.Range("G79").Value = 5.5407960170913082
$ws.Range("H79").Value = 51.588250408325628
$ws.Range("G80").Value = 5.5635885158603031
$ws.Range("H80").Value = 51.628362152281937
$ws.Range("G81").Value = 5.563791314933872
$ws.Range("H81").Value = 51.623856203544122
$ws.Range("G82").Value = 5.5522097629289986
$ws.Range("H82").Value = 51.62398882963442
$ws.Range("G83").Value = 5.5537412260984853
$ws.Range("H83").Value = 51.606170025655658
$ws.Range("G84").Value = 5.5216310171119023
$ws.Range("H84").Value = 51.597283835410742
$ws.Range("G85").Value = 5.5149267287937276
$ws.Range("H85").Value = 51.610691841517912
$ws.Range("G86").Value = 5.4973898688597549
$ws.Range("H86").Value = 51.604374139185452
$ws.Range("G87").Value = 5.5673201305593922
$ws.Range("H87").Value = 51.603904966114108
$ws.Range("G88").Value = 5.5537527260778701
$ws.Range("H88").Value = 51.627201716280403
$ws.Range("G89").Value = 5.5402225618773953
$ws.Range("H89").Value = 51.61365786673359
$ws.Range("G90").Value = 5.545246278073761
$ws.Range("H90").Value = 51.622914108383647
$ws.Range("G91").Value = 5.5574583971830949
$ws.Range("H91").Value = 51.620096358880247
$ws.Range("G92").Value = 5.5399848298155661
$ws.Range("H92").Value = 51.606960663111003
$ws.Range("G93").Value = 5.5428949868554964
$ws.Range("H93").Value = 51.610557282077721
$ws.Range("G94").Value = 5.5805300583021147
$ws.Range("H94").Value = 51.633868793760342
$ws.Range("G95").Value = 5.5346674511126661
$ws.Range("H95").Value = 51.616920391220113
$ws.Range("G96").Value = 5.5618473927866363
$ws.Range("H96").Value = 51.613462744034528
$ws.Range("G97").Value = 5.545696916319014
$ws.Range("H97").Value = 51.617579893262089
$ws.Range("G98").Value = 5.5497305775696804
$ws.Range("H98").Value = 51.612453521387607
$ws.Range("G99").Value = 5.5515756654829369
$ws.Range("H99").Value = 51.609402825875478
$ws.Range("G100").Value = 5.5378596006976739
$ws.Range("H100").Value = 51.621359439959321
$ws.Range("G101").Value = 5.5407346193091849
$ws.Range("H101").Value = 51.592746186588371

# Update active cell selection
$ws.Range("P26").Select()
